$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questions")
$r = $ws.Range("E7")
$r.Font.Name = "Segoe UI"
